$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2755.3635
$ws.Range("I28").Value = 539.125
$ws.Range("J28").Value = 8665.333000000001
$ws.Range("K28").Value = 539.125
$ws.Range("L28").Value = 8665.333000000001
$ws.Range("M28").Value = -54.125
$ws.Range("N28").Value = -9635.333000000001

$ws.Range("H50").Value = 198
$ws.Range("J50").Value = 198
$ws.Range("L50").Value = 594
$ws.Range("N50").Value = -1544

$ws.Range("H98").Value = 3224.7368
$ws.Range("I98").Value = 3554.7058
$ws.Range("K98").Value = 3554.7058
$ws.Range("M98").Value = -2056.7058

$ws.Range("H111").Value = 5000
$ws.Range("I111").Value = 5000
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 15000
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -11933
$ws.Range("N111").Value = ""

$ws.Range("H122").Value = 3224.7368
$ws.Range("I122").Value = 3554.7058
$ws.Range("K122").Value = 10664.1174
$ws.Range("M122").Value = -8214.117400000001

$ws.Range("H137").Value = 1223269.8
$ws.Range("I137").Value = 5556367.5
$ws.Range("J137").Value = 4586.0625
$ws.Range("K137").Value = 16669102.5
$ws.Range("L137").Value = 13758.1875
$ws.Range("M137").Value = -16666552.5
$ws.Range("N137").Value = -18858.1875

$ws.Range("H138").Value = 3118.24
$ws.Range("J138").Value = 2734.8
$ws.Range("L138").Value = 8204.400000000001
$ws.Range("N138").Value = -18484.4

$ws.Range("H141").Value = 2935.6
$ws.Range("I141").Value = 2702.2307
$ws.Range("K141").Value = 8106.6921
$ws.Range("M141").Value = -2926.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3515.6038
$ws.Range("I32").Value = 3167.068
$ws.Range("J32").Value = 5219.5557
$ws.Range("K32").Value = 3167.068
$ws.Range("L32").Value = 5219.5557
$ws.Range("M32").Value = -2880.068
$ws.Range("N32").Value = -5793.5557

$ws.Range("H61").Value = 3588.75
$ws.Range("I61").Value = 900
$ws.Range("J61").Value = 5202
$ws.Range("K61").Value = 900
$ws.Range("L61").Value = 5202
$ws.Range("M61").Value = -688
$ws.Range("N61").Value = -5626

$ws.Range("H74").Value = 224081.36
$ws.Range("I74").Value = 348270
$ws.Range("J74").Value = 3301.5557
$ws.Range("K74").Value = 348270
$ws.Range("L74").Value = 3301.5557
$ws.Range("M74").Value = -347396
$ws.Range("N74").Value = -5049.5557

$ws.Range("H77").Value = 224081.36
$ws.Range("I77").Value = 348270
$ws.Range("J77").Value = 3301.5557
$ws.Range("K77").Value = 1741350
$ws.Range("L77").Value = 16507.7785
$ws.Range("M77").Value = -1736982
$ws.Range("N77").Value = -25243.7785

$ws.Range("H132").Value = 1461.069
$ws.Range("I132").Value = 723.7917
$ws.Range("K132").Value = 2171.3751
$ws.Range("M132").Value = 358.6248999999998

$ws.Range("H136").Value = 3588.75
$ws.Range("I136").Value = 900
$ws.Range("J136").Value = 5202
$ws.Range("K136").Value = 2700
$ws.Range("L136").Value = 15606
$ws.Range("M136").Value = -150
$ws.Range("N136").Value = -20706

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1435.0344
$ws.Range("I16").Value = 1287.4546
$ws.Range("K16").Value = 1287.4546
$ws.Range("M16").Value = -1000.4546

$ws.Range("H31").Value = 6072.696
$ws.Range("I31").Value = 4383.778
$ws.Range("J31").Value = 7158.4287
$ws.Range("K31").Value = 4383.778
$ws.Range("L31").Value = 7158.4287
$ws.Range("M31").Value = -4088.778
$ws.Range("N31").Value = -7748.4287

$ws.Range("H34").Value = 6072.696
$ws.Range("I34").Value = 4383.778
$ws.Range("J34").Value = 7158.4287
$ws.Range("K34").Value = 4383.778
$ws.Range("L34").Value = 7158.4287
$ws.Range("M34").Value = -4181.778
$ws.Range("N34").Value = -7562.4287

$ws.Range("H87").Value = 59999.5
$ws.Range("J87").Value = 59999.5
$ws.Range("L87").Value = 59999.5
$ws.Range("N87").Value = -62371.5

$ws.Range("H90").Value = 59999.5
$ws.Range("J90").Value = 59999.5
$ws.Range("L90").Value = 179998.5
$ws.Range("N90").Value = -191854.5

$ws.Range("H99").Value = 3787.9
$ws.Range("I99").Value = 3680
$ws.Range("K99").Value = 3680
$ws.Range("M99").Value = -2182

$ws.Range("H103").Value = 4677.2856
$ws.Range("I103").Value = 4677.2856
$ws.Range("K103").Value = 4677.2856
$ws.Range("M103").Value = -3505.2856

$ws.Range("H113").Value = 1435.0344
$ws.Range("I113").Value = 1287.4546
$ws.Range("K113").Value = 1287.4546
$ws.Range("M113").Value = 882.5454

$ws.Range("H126").Value = 3787.9
$ws.Range("I126").Value = 3680
$ws.Range("K126").Value = 11040
$ws.Range("M126").Value = -8570

$ws.Range("H132").Value = 2563.1538
$ws.Range("I132").Value = 2565.125
$ws.Range("J132").Value = 2560
$ws.Range("K132").Value = 7695.375
$ws.Range("L132").Value = 7680
$ws.Range("M132").Value = -5165.375
$ws.Range("N132").Value = -12740

$ws.Range("H134").Value = 2527.5806
$ws.Range("I134").Value = 2309.4443
$ws.Range("K134").Value = 6928.3329
$ws.Range("M134").Value = -4393.3329

$ws.Range("H141").Value = 477476.8
$ws.Range("J141").Value = 477476.8
$ws.Range("L141").Value = 477476.8
$ws.Range("N141").Value = -487836.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2423.7273
$ws.Range("I3").Value = 1522.5714
$ws.Range("J3").Value = 4000.75
$ws.Range("K3").Value = 4567.7142
$ws.Range("L3").Value = 12002.25
$ws.Range("M3").Value = -4455.7142
$ws.Range("N3").Value = -12226.25

$ws.Range("H5").Value = 1513.1875
$ws.Range("I5").Value = 594.5
$ws.Range("J5").Value = 1819.4166
$ws.Range("K5").Value = 1783.5
$ws.Range("L5").Value = 5458.2498
$ws.Range("M5").Value = -1671.5
$ws.Range("N5").Value = -5682.2498

$ws.Range("H76").Value = 7429.6665
$ws.Range("J76").Value = 7444.5
$ws.Range("L76").Value = 22333.5
$ws.Range("N76").Value = -23099.5

$ws.Range("H79").Value = 7429.6665
$ws.Range("J79").Value = 7444.5
$ws.Range("L79").Value = 22333.5
$ws.Range("N79").Value = -24985.5

$ws.Range("H135").Value = 1513.1875
$ws.Range("I135").Value = 594.5
$ws.Range("J135").Value = 1819.4166
$ws.Range("K135").Value = 5350.5
$ws.Range("L135").Value = 16374.7494
$ws.Range("M135").Value = -2815.5
$ws.Range("N135").Value = -21444.7494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 853
$ws.Range("I16").Value = 842.5
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 842.5
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -672.5
$ws.Range("N16").Value = -1340

$ws.Range("H22").Value = 1049.8823
$ws.Range("I22").Value = 1271.7
$ws.Range("J22").Value = 733
$ws.Range("K22").Value = 1271.7
$ws.Range("L22").Value = 733
$ws.Range("M22").Value = -976.7
$ws.Range("N22").Value = -1323

$ws.Range("H27").Value = 1049.8823
$ws.Range("I27").Value = 1271.7
$ws.Range("J27").Value = 733
$ws.Range("K27").Value = 1271.7
$ws.Range("L27").Value = 733
$ws.Range("M27").Value = -1164.7
$ws.Range("N27").Value = -947

$ws.Range("H46").Value = 2995.4546
$ws.Range("I46").Value = 2745
$ws.Range("K46").Value = 2745
$ws.Range("M46").Value = -2557

$ws.Range("H61").Value = 7240.0527
$ws.Range("I61").Value = 1617.0667
$ws.Range("K61").Value = 1617.0667
$ws.Range("M61").Value = -1415.0667

$ws.Range("H113").Value = 7240.0527
$ws.Range("I113").Value = 1617.0667
$ws.Range("K113").Value = 1617.0667
$ws.Range("M113").Value = 552.9332999999999

$ws.Range("H127").Value = 63499.5
$ws.Range("J127").Value = 63499.5
$ws.Range("L127").Value = 63499.5
$ws.Range("N127").Value = -73419.5

$ws.Range("H132").Value = 3205.9048
$ws.Range("I132").Value = 1746.75
$ws.Range("J132").Value = 4103.846
$ws.Range("K132").Value = 5240.25
$ws.Range("L132").Value = 12311.538
$ws.Range("M132").Value = -2710.25
$ws.Range("N132").Value = -17371.538

$ws.Range("H134").Value = 103977.664
$ws.Range("J134").Value = 103977.664
$ws.Range("L134").Value = 103977.664
$ws.Range("N134").Value = -114117.664

$ws.Range("H137").Value = 57940.41
$ws.Range("J137").Value = 57940.41
$ws.Range("L137").Value = 57940.41
$ws.Range("N137").Value = -68140.41

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 15712.833
$ws.Range("J56").Value = 15712.833
$ws.Range("L56").Value = 15712.833
$ws.Range("N56").Value = -17140.833

$ws.Range("H126").Value = 2401.2727
$ws.Range("I126").Value = 1502.3334
$ws.Range("J126").Value = 3480
$ws.Range("K126").Value = 4507.0002
$ws.Range("L126").Value = 10440
$ws.Range("M126").Value = -2037.0002
$ws.Range("N126").Value = -15380

$ws.Range("H140").Value = 102062.5
$ws.Range("J140").Value = 98015.71000000001
$ws.Range("L140").Value = 98015.71000000001
$ws.Range("N140").Value = -108375.71
